# Estado de Cuenta - agrega nuevo trabajador (ANDRES MEDINA AGUILAR) con 16
# periodos de mora, conserva los registros previos (que se recorren hacia
# abajo) y actualiza los totales del encabezado.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Insertar 16 filas en blanco justo encima de los datos existentes
#    (fila 16). Esto recorre los registros actuales (filas 16-27) hacia
#    las filas 32-43, y el bloque de firmas (filas 32-33) hacia 48-49.
# ---------------------------------------------------------------------
$ws.Range("16:31").Insert()

# Copiar el formato de una fila de datos "normal" (ahora en la fila 42,
# antes fila 26) hacia las filas recien insertadas, para que luzcan
# igual que el resto de la tabla.
$ws.Range("B42:J42").Copy()
$ws.Range("B16:J31").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2) Volcar la tabla completa (filas 16 a 43): primero los 16 periodos
#    de mora del nuevo trabajador ANDRES MEDINA AGUILAR (CC 73185464)
#    y luego los registros previos, que quedan reordenados tal como en
#    la nueva base de datos.
# ---------------------------------------------------------------------
$tabla = @(
  @("CC","73185464","ANDRES MEDINA AGUILAR","2209",30430,1082828),
  @("CC","73185464","ANDRES MEDINA AGUILAR","2208",35112,1082828),
  @("CC","73185464","ANDRES MEDINA AGUILAR","2207",35112,1082828),
  @("CC","73185464","ANDRES MEDINA AGUILAR","2206",35112,1082828),
  @("CC","73185464","ANDRES MEDINA AGUILAR","2205",35112,1082828),
  @("CC","73185464","ANDRES MEDINA AGUILAR","2204",35112,1082828),
  @("CC","73185464","ANDRES MEDINA AGUILAR","2203",35112,1082828),
  @("CC","73185464","ANDRES MEDINA AGUILAR","2202",35112,1082828),
  @("CC","73185464","ANDRES MEDINA AGUILAR","2110",35112,1082828),
  @("CC","73185464","ANDRES MEDINA AGUILAR","2109",35112,1082828),
  @("CC","73185464","ANDRES MEDINA AGUILAR","2108",36341,1082828),
  @("CC","73185464","ANDRES MEDINA AGUILAR","2107",36341,1082828),
  @("CC","73185464","ANDRES MEDINA AGUILAR","2106",36341,1082828),
  @("CC","73185464","ANDRES MEDINA AGUILAR","2105",36341,1082828),
  @("CC","73185464","ANDRES MEDINA AGUILAR","2104",36341,1082828),
  @("CC","73185464","ANDRES MEDINA AGUILAR","2103",36341,1082828),
  @("CC","73572087","EDUARDO HERNANDEZ CABARCAS","2202",40000,1000000),
  @("CC","40987503","GLORIA VICTORIA VANEGAS REYES","2203",40000,1000000),
  @("CC","40987503","GLORIA VICTORIA VANEGAS REYES","2202",40000,1000000),
  @("CC","94373000","FRANK YAIR CHAUX AVILA","2209",34666,1000000),
  @("CC","1043964778","NELSON ENRIQUE PACHECO BOHORQUEZ","2209",34666,1000000),
  @("CC","1043964778","NELSON ENRIQUE PACHECO BOHORQUEZ","2208",40000,1000000),
  @("CC","1043964778","NELSON ENRIQUE PACHECO BOHORQUEZ","2207",40000,1000000),
  @("CC","1043964778","NELSON ENRIQUE PACHECO BOHORQUEZ","2206",40000,1000000),
  @("CC","1043964778","NELSON ENRIQUE PACHECO BOHORQUEZ","2205",40000,1000000),
  @("CC","1043964778","NELSON ENRIQUE PACHECO BOHORQUEZ","2204",40000,1000000),
  @("CC","1043964778","NELSON ENRIQUE PACHECO BOHORQUEZ","2203",40000,1000000),
  @("CC","1043964778","NELSON ENRIQUE PACHECO BOHORQUEZ","2202",40000,1000000)
)

for ($i = 0; $i -lt $tabla.Count; $i++) {
    $fila = 16 + $i
    $reg = $tabla[$i]
    $ws.Cells.Item($fila, 2).Value = $reg[0]
    $ws.Cells.Item($fila, 3).Value = $reg[1]
    $ws.Cells.Item($fila, 4).Value = $reg[2]
    $ws.Cells.Item($fila, 5).Value = $reg[3]
    $ws.Cells.Item($fila, 6).Value = $reg[4]
    $ws.Cells.Item($fila, 7).Value = $reg[5]
}

# ---------------------------------------------------------------------
# 3) Actualizar los totales de la cabecera del estado de cuenta.
# ---------------------------------------------------------------------
$ws.Range("E11").Value = 1033816   # Valor mora total
$ws.Range("C13").Value = 5         # Cant. Trabajadores
$ws.Range("F13").Value = 16        # Cant. Periodos
